$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the left (old B:G -> new A:F), dropping the
# old (empty-header) column A. This also naturally drops the bold/border
# style that used to live on column A's data cells, since the incoming
# column (old B) never had that style.
$ws.Range("A1").EntireColumn.Delete()

# The sheet only has rows through 17; the well now needs three more survey
# stations. Duplicate the structure of the last existing row (which still
# carries blank x/y/z cells) down into rows 18-20 before overwriting the
# md/inclination/azimuth values, so the blank x/y/z cells keep existing
# (rather than simply being absent) on the new rows too.
$ws.Range("A17:F17").Copy($ws.Range("A18:F18"))
$ws.Range("A17:F17").Copy($ws.Range("A19:F19"))
$ws.Range("A17:F17").Copy($ws.Range("A20:F20"))

# Row 2 (md=1 / straight start of the well)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1

# Row 3 (x/y/z stay blank - they already are after the column shift)
$ws.Range("A3").Value = 1001
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0

# Remaining survey rows (md, inclination, azimuth); x/y/z left blank.
$surveyRows = @(
    @{ Row = 4;  Md = 1031; Inc = 5.399999999999999;  Az = 123 },
    @{ Row = 5;  Md = 1061; Inc = 10.8;                Az = 123 },
    @{ Row = 6;  Md = 1091; Inc = 16.2;                Az = 123 },
    @{ Row = 7;  Md = 1121; Inc = 21.6;                Az = 123 },
    @{ Row = 8;  Md = 1151; Inc = 27;                  Az = 123 },
    @{ Row = 9;  Md = 1181; Inc = 32.4;                Az = 123 },
    @{ Row = 10; Md = 1211; Inc = 37.8;                Az = 123 },
    @{ Row = 11; Md = 1241; Inc = 43.2;                Az = 123 },
    @{ Row = 12; Md = 1271; Inc = 48.59999999999999;  Az = 123 },
    @{ Row = 13; Md = 1301; Inc = 53.99999999999999;  Az = 123 },
    @{ Row = 14; Md = 1331; Inc = 59.39999999999999;  Az = 123 },
    @{ Row = 15; Md = 1361; Inc = 64.8;                Az = 123 },
    @{ Row = 16; Md = 1391; Inc = 70.19999999999999;  Az = 123 },
    @{ Row = 17; Md = 1421; Inc = 75.59999999999999;  Az = 123 },
    @{ Row = 18; Md = 1451; Inc = 80.99999999999999;  Az = 123 },
    @{ Row = 19; Md = 1481; Inc = 86.39999999999999;  Az = 123 },
    @{ Row = 20; Md = 1881; Inc = 90;                  Az = 123 }
)

foreach ($r in $surveyRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Md
    $ws.Range("B$row").Value = $r.Inc
    $ws.Range("C$row").Value = $r.Az
}
